$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 4.670199999999999
$ws.Range("A12").Value = -22.81520000000002
$ws.Range("B12").Value = 5.598700000000003
$ws.Range("B14").Value = 8.758400000000005
$ws.Range("B22").Value = 4.731800000000004
